$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the weekly salary figures (B8:B11) ---
$ws.Range("B8").Value = 90
$ws.Range("B9").Value = 90
$ws.Range("B10").Value = 110
$ws.Range("B11").Value = 110

# --- Update "Tasks completed" / "Tasks to complete" table (rows 19-21) ---
# Rows 19-21 (both columns) get the smaller 10pt font used for the new entries.
$ws.Range("A19:B21").Font.Size = 10

$ws.Range("A19").Value = "Finished the manual design of the frontend."
$ws.Range("B19").Value = "Automate our Hi-Fi prototype."

$ws.Range("A20").Value = "Worked on the backend."
$ws.Range("B20").ClearContents()

$ws.Range("A21").Value = "Prepared the presentation of the Hi-Fi prototype."
$ws.Range("B21").ClearContents()

# Column B of the remaining blank rows (22-25) also switches to the 10pt font,
# while column A keeps its original formatting.
$ws.Range("B22:B25").Font.Size = 10

# --- Selection cursor moves to B21 ---
$ws.Range("B21").Select()

# --- Page setup: A4, portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
